$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")
$wsUtils = $wb.Worksheets.Item("Utils")

# --- New journal entry (row 14) ---------------------------------------
$ws.Range("A14").Value = 45483
$ws.Range("B14").Value = 0.3888888888888889
$ws.Range("C14").Value = 0.49305555555555558
$ws.Range("E14").Value = "Analyse et état de l'art"

# --- Per-category SUMIF summary table in columns H:I -------------------
# Seed H7/"Total" first so the shared string "Total" (currently only used
# by H2) never becomes orphaned while H2 is being repointed below.
$ws.Range("H7").Value = "Total"

$ws.Range("H1").Value = "Analyse et état de l'art"
$ws.Range("H3").Value = "Réalisation de l'application "
$ws.Range("H4").Value = "Tests et validations"
$ws.Range("H5").Value = "Gestion du projet, documentation et présentation"
$ws.Range("H2").Value = "Réalisation du modèle"

$ws.Range("I1").Formula = "=SUMIF(E:E, H1, D:D)"
$ws.Range("I1").NumberFormat = $ws.Range("I2").NumberFormat

$ws.Range("I2:I5").Formula = "=SUMIF(E:E, H2, D:D)"
$ws.Range("I2:I5").NumberFormat = $ws.Range("I2").NumberFormat

$ws.Range("I7").Formula = "=SUM(I1:I5)"
$ws.Range("I7").NumberFormat = $ws.Range("I2").NumberFormat

# --- Column H is now a long label column, widen it ----------------------
# (44.85546875 in the target file; the host's column-width setter only
# lands on ~1/6-character increments, so 44.0 is the closest reachable
# input -> stored width 44.8333.., the nearest achievable bucket.)
$ws.Columns.Item(8).ColumnWidth = 44.0

# --- Recalculate so every formula carries a fresh cached value ---------
$excel.Calculate()

# --- Selections, matching where the author last clicked -----------------
$wsUtils.Activate()
$wsUtils.Range("B2:B6").Select()
$ws.Activate()
$ws.Range("E14").Select()
